$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sample/test order rows (4 through 15) contained dummy data used while
# testing the newly-added "Received At" date column. Clear that leftover
# sample data, leaving the existing cell formatting (styles) untouched so the
# template rows remain ready to be filled by the Kafka consumer app.
for ($row = 4; $row -le 15; $row++) {
    $rowRange = $ws.Range("B" + $row + ":K" + $row)
    $rowRange.ClearContents()
}

# Match the author's final cursor position in the sheet.
$ws.Range("J7").Select()
